$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "55.993.70"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.53%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.500.30"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "486.61"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +5.20%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "145.20"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +11.64%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.509"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.95%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.519.02"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.58%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "5.65"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.30%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0973"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.83%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.331"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.81%  "
$ws.Range("E13").Value = "  +1.07%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.928.62"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.38%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "55.987.65"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.73%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.03"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +6.97%  "
$ws.Range("E17").Value = "  +3.39%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.505.50"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.03%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.46"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.75%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.25"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +9.33%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "319.88"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.33%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.80"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +7.78%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "58.36"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.26%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.410"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +6.86%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.166"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +8.34%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.602.70"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.81%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.54"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +5.03%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0₃0784"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +10.00%  "
$ws.Range("E31").Value = "  +0.27%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "148.21"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "18.33"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.85%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.50"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +6.97%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.22"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.38%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.15"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +9.21%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.71"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +5.03%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.866"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +8.53%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "34.27"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.26%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.54"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.65%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.616"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("E42").Value = "  +0.03%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0554"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.68%  "
$ws.Range("E44").Value = "  +6.11%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.81"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +10.53%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "261.67"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +20.44%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "10.16"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  +3.53%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0904"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.32%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.922.49"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.42%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "17.66"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +6.81%  "
